$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.148.12'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.020.62'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = '226.98'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').Value = '0.608'
$ws.Range('E6').Value = '  -2.13%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '55.06'
$ws.Range('E8').Value = '  -3.69%  '
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('D10').Value = '0.0789'
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('D12').Value = '2.319.14'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').Value = '14.32'
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('D14').Value = '20.47'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '0.743'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('D16').Value = '5.15'
$ws.Range('E16').Value = '  -1.96%  '
$ws.Range('D17').Value = '2.022.25'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = '37.065.50'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '6.18'
$ws.Range('E19').Value = '  +3.36%  '
$ws.Range('D20').Value = '69.00'
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').Value = '226.11'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('E25').Value = '  -5.35%  '
$ws.Range('D26').Value = '165.55'
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').Value = '9.21'
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('E28').Value = '  -2.79%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '18.76'
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('D33').Value = '0.0619'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').Value = '4.43'
$ws.Range('E34').Value = '  -3.20%  '
$ws.Range('E35').Value = '  -4.64%  '
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('D39').Value = '5.42'
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('E40').Value = '  -4.12%  '
$ws.Range('D41').Value = '1.483.72'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D44').Value = '0.0926'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D47').Value = '7.30'
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').Value = '2.92'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = '2.208.28'
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('D51').Value = '44.44'
$ws.Range('E51').Value = '  -1.94%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '95.42'
$ws.Range('E42').Value = '  -2.73%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '16.64'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.76'
$ws.Range('E45').Value = '  -4.58%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '1.14'
$ws.Range('E46').Value = '  -4.32%  '
